$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.264.06"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.237.69"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'306.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").Value = "'93.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.72%  "
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").Value = "'34.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.62%  "
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.41%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "2.332.38"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").Value = "'0.824"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").Value = "'13.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("D17").Value = "43.947.63"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "0.0₃0961"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "'11.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.31%  "
$ws.Range("D21").Value = "'65.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").Value = "'236.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "'1.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'39.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").Value = "'9.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.17%  "
$ws.Range("D29").Value = "'19.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "'5.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("D31").Value = "'151.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("D32").Value = "'0.0789"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.65%  "
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("D34").Value = "'3.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.76%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("D37").Value = "'1.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.97%  "
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").Value = "'14.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("D40").Value = "'3.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "1.704.09"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "'82.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("D47").Value = "'98.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("D48").Value = "'1.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("D49").Value = "'54.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").Value = "'8.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").Value = "'66.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.73%  "
